$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Remove rows 4-17, shifting the sheet dimension down to A1:F3
$ws.Range("A4:F17").EntireRow.Delete()

# Row 2 becomes the "OTROS" group with all-zero figures
$ws.Cells.Item(2, 2).Value = "OTROS"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0

# Row 3 becomes the TOTAL row: drop the ASESOR label, relabel GRUPO -> TOTAL
$ws.Cells.Item(3, 1).ClearContents()
$ws.Cells.Item(3, 2).Value = "TOTAL"
$ws.Cells.Item(3, 2).HorizontalAlignment = -4152
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0

# Narrower column widths for the now-shorter table
$ws.Columns.Item(2).ColumnWidth = 6.166666666666667
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 16.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668
